$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: Volume number 49 -> 50
$ws.Range("A8").Characters(21, 2).Text = "50"

# Header: report week dates
$ws.Range("C9").Characters(27, 9).Text = "12/12/2022"
$ws.Range("C9").Characters(48, 10).Text = "12/18/2022"

# --- Table data updates rows 15-30 ---
# Row 15
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "0"
$ws.Range("D15").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("L15").Value = 50
$ws.Range("N15").Value = -16

# Row 16
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "0"
$ws.Range("D15").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = -22.222222222222
$ws.Range("J16").Value = 129
$ws.Range("K16").Value = 17.054263565891
$ws.Range("L16").Value = 17.054263565891
$ws.Range("N16").Value = -88.739746457867

# Row 17
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 200
$ws.Range("F17").Value = 7
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = -46.153846153846
$ws.Range("I17").Value = 173
$ws.Range("J17").Value = 145
$ws.Range("K17").Value = 19.310344827586
$ws.Range("L17").Value = 39.516129032258
$ws.Range("M17").Value = 31.060606060606
$ws.Range("N17").Value = -63.807531380753

# Row 18
$ws.Range("D18").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = -72.727272727272
$ws.Range("I18").Value = 199
$ws.Range("J18").Value = 193
$ws.Range("K18").Value = 3.108808290155
$ws.Range("L18").Value = -23.461538461538
$ws.Range("M18").Value = 9.944751381215
$ws.Range("N18").Value = -90.432692307692

# Row 19
$ws.Range("C19").Value = 53
$ws.Range("D19").Value = 45
$ws.Range("E19").Value = 17.777777777777
$ws.Range("F19").Value = 202
$ws.Range("G19").Value = 199
$ws.Range("H19").Value = 1.507537688442
$ws.Range("I19").Value = 1947
$ws.Range("J19").Value = 1272
$ws.Range("K19").Value = 53.066037735849
$ws.Range("L19").Value = 86.673058485139
$ws.Range("M19").Value = 12.673611111111
$ws.Range("N19").Value = -73.470500068129

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0"
$ws.Range("D15").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "***.*"
$ws.Range("H15").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = -37.5
$ws.Range("I20").Value = 114
$ws.Range("K20").Value = 58.333333333333
$ws.Range("L20").Value = 65.217391304347
$ws.Range("M20").Value = 153.333333333333
$ws.Range("N20").Value = -76.446280991735

# Row 21
$ws.Range("C21").Value = 58
$ws.Range("D21").Value = 52
$ws.Range("E21").Value = 11.538461538461
$ws.Range("F21").Value = 225
$ws.Range("G21").Value = 240
$ws.Range("H21").Value = -6.25
$ws.Range("I21").Value = 2606
$ws.Range("J21").Value = 1832
$ws.Range("K21").Value = 42.248908296943
$ws.Range("L21").Value = 58.612294583079
$ws.Range("M21").Value = 16.495306213679
$ws.Range("N21").Value = -77.836366729035

# Row 22
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("D15").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F22").Value = 6
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "0"
$ws.Range("D15").Copy()
$ws.Range("G22").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("H22").NumberFormat = "@"
$ws.Range("H22").Value = "***.*"
$ws.Range("H15").Copy()
$ws.Range("H22").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("L22").Value = -9.459459459459
$ws.Range("M22").Value = 3.076923076923

# Row 23
$ws.Range("F23").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C23").Value = 1
$ws.Range("I23").Value = 4
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 100
$ws.Range("M23").Value = -33.333333333333

# Row 24
$ws.Range("C24").Value = 38
$ws.Range("D24").Value = 36
$ws.Range("E24").Value = 5.555555555555
$ws.Range("F24").Value = 226
$ws.Range("G24").Value = 166
$ws.Range("H24").Value = 36.144578313253
$ws.Range("I24").Value = 2676
$ws.Range("J24").Value = 1915
$ws.Range("K24").Value = 39.738903394255
$ws.Range("L24").Value = 92.934390771449
$ws.Range("M24").Value = 41.213720316622

# Row 25
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = -37.5
$ws.Range("F25").Value = 38
$ws.Range("G25").Value = 35
$ws.Range("H25").Value = 8.571428571428
$ws.Range("I25").Value = 549
$ws.Range("J25").Value = 439
$ws.Range("K25").Value = 25.056947608200
$ws.Range("L25").Value = 81.788079470198
$ws.Range("M25").Value = 23.094170403587

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0"
$ws.Range("D15").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "***.*"
$ws.Range("H15").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F26").Value = 4
$ws.Range("H26").Value = 300
$ws.Range("I26").Value = 37
$ws.Range("K26").Value = 19.354838709677
$ws.Range("L26").Value = 76.190476190476

# Row 27
$ws.Range("C27").Value = 1
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("D15").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "***.*"
$ws.Range("H15").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F27").Value = 9
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 80
$ws.Range("I27").Value = 100
$ws.Range("K27").Value = 20.481927710843
$ws.Range("L27").Value = 72.413793103448

# Row 30
$ws.Range("G30").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D30").Value = 1
$ws.Range("H30").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E30").Value = -100
$ws.Range("G30").Value = 1
$ws.Range("J30").Value = 18
$ws.Range("K30").Value = -16.666666666666
